$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The report lost two product lines ("FORFLOZIN 10MG 30 F.C. TABS" on the
# original row 6, and "TAVONIZA 20 MG 20 F.C.TABS." on the original row 12).
# Delete the lower row first so the upper row's index doesn't shift before
# it is removed.
$ws.Rows("12:12").Delete()
$ws.Rows("6:6").Delete()

# The "م" serial-number column (A) holds plain literal numbers, not a
# formula, so it needs to be renumbered 1..11 after the two rows were
# removed (rows now run 4..14).
for ($i = 0; $i -lt 11; $i++) {
    $ws.Cells.Item(4 + $i, 1).Value = $i + 1
}

# The grand-total cell (now on row 15 after the shift) is a literal sum,
# not a formula, so update it to drop the two deleted rows' prices
# (102 + 99 = 201 -> 610.11 - 201 = 409.11).
$ws.Range("K15").Value = 409.11

# The per-row heights are fixed formatting (not auto-carried with the
# deleted rows), so restore rows 4..14 to their original heights.
$heights = @{4=24.75; 5=25.5; 6=24.75; 7=25.5; 8=25.5; 9=24.75; 10=25.5; 11=24.75; 12=25.5; 13=25.5; 14=24.75}
foreach ($r in $heights.Keys) {
    $ws.Rows("$($r):$($r)").RowHeight = $heights[$r]
}
